$d = $word.ActiveDocument

# --- Step 1: rewrite the paragraphs that survive, in place -------------
# (Done before any deletions so the original 1-based Paragraphs.Item
# indices still line up with the pre-edit document.)

# Block 1: validate_password -> is_valid_password body
$d.Paragraphs.Item(3).Range.Text  = 'def is_valid_password(password):'
$d.Paragraphs.Item(4).Range.Text  = '    if len(password) < 8:'
$d.Paragraphs.Item(5).Range.Text  = '        return False, ''Password must be at least 8 characters long.'' '
$d.Paragraphs.Item(6).Range.Text  = '    if not re.search(r''\d'', password):'
$d.Paragraphs.Item(7).Range.Text  = '        return False, ''Password must include at least one number.'''
$d.Paragraphs.Item(8).Range.Text  = '    if not re.search(r''[!@#$%^&*(),.?":{}|<>]'', password):'
$d.Paragraphs.Item(9).Range.Text  = '        return False, ''Password must include at least one special character.'''
$d.Paragraphs.Item(10).Range.Text = '    return True, ''Password is valid.'''

# Block 2: test_passwords() body swapped for a loop-driven version
$d.Paragraphs.Item(17).Range.Text = '    test_cases = [''abc123'', ''abc12345'', ''abc123@'', ''mypassword1'', ''Pass123!'', ''12345678'', ''MyPass@'']'
$d.Paragraphs.Item(18).Range.Text = '    for pwd in test_cases:'
$d.Paragraphs.Item(19).Range.Text = '        print(f''Password: {pwd} - {message}'')'

# --- Step 2: drop the paragraphs that no longer have a home -------------
# Old paragraphs 11-14 (leftover stray return/blank lines after the
# special-character check) and old paragraphs 20-47 (the ten individual
# test1..test10 call/assert pairs, replaced by the test_cases loop).
# Delete from the highest index down so lower indices stay valid.
$toDelete = @(47,46,45,44,43,42,41,40,39,38,37,36,35,34,33,32,31,30,29,28,27,26,25,24,23,22,21,20,14,13,12,11)
foreach ($i in $toDelete) {
    $d.Paragraphs.Item($i).Range.Delete()
}
